$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58 - this pushes the existing rows 58..86
# down to 59..87 (matching the diff, which shows every row from 58
# downward taking on the values that previously belonged to the row
# above it, and a brand-new row 87 appearing with what used to be
# row 86's data).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly price record.
$ws.Range("A58").Value = 9
$ws.Range("B58").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 45176
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = 100112010
$ws.Range("G58").Value = "Achicoria"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 90
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = 7000
$ws.Range("N58").Value = '$/caja 16 unidades'
$ws.Range("O58").Value = "Provincia de Quillota"
$ws.Range("P58").Value = 438
$ws.Range("Q58").Value = 16
$ws.Range("R58").Value = "Hortaliza"
